$d = $word.ActiveDocument

# Change 1: merge "который ждет когда" sequence (remove grammar proofErr split)
$d.Content.Find.Execute("который ждет когда", $true, $false, $false, $false, $false, $true, 1, $false, "который ждет когда", 2) | Out-Null

# Change 2: fix misspelling "совреминем" -> "со временем"
$d.Content.Find.Execute("совреминем", $true, $false, $false, $false, $false, $true, 1, $false, "со временем", 2) | Out-Null
